# Add the "28. 9. 2021" data wave as a new last column on both sheets,
# and bump the "aktualizace" (last-updated) date in the footer/title rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "data" (percentages) - new column AH, header "28. 9. 2021"
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")

# Copy the header formatting from the previous last header cell (AG1) so the
# new header cell (AH1) looks the same (bold, centered, bordered).
$wsData.Range("AG1").Copy()
$wsData.Range("AH1").PasteSpecial(-4122)  # xlPasteFormats
$wsData.Range("AH1").Value = "28. 9. 2021"

$dataValues = @{
    2  = 0.09
    3  = 0.07
    4  = 0.11
    5  = 0.05
    6  = 0.08
    7  = 0.12
    8  = 0.11
    9  = 0.15
    10 = 0.11
    11 = 0.08
    12 = 0.07
    13 = 0.09
    14 = 0.22
    15 = 0.11
    16 = 0.07
    17 = 0.14
    18 = 0.11
    19 = 0.06
    20 = 0.08
    21 = 0.06
    22 = 0.06
    23 = 0.17
}

foreach ($row in $dataValues.Keys) {
    $wsData.Range("AH$row").Value = $dataValues[$row]
}

# Bump the "aktualizace" date mentioned in the title row at the bottom of the sheet.
$wsData.Range("A24").Value = "Život během pandemie, Duševní zdraví, % respondentů celkově a ve skupinách, aktualizace 6. 10. 2021"

# ---------------------------------------------------------------------
# Sheet "pocetR" (sample sizes) - new column AG, header "28. 9. 2021"
# ---------------------------------------------------------------------
$wsPocet = $wb.Worksheets.Item("pocetR")

$wsPocet.Range("AF1").Copy()
$wsPocet.Range("AG1").PasteSpecial(-4122)  # xlPasteFormats
$wsPocet.Range("AG1").Value = "28. 9. 2021"

$pocetValues = @{
    2  = 1855
    3  = 897
    4  = 958
    5  = 230
    6  = 667
    7  = 283
    8  = 675
    9  = 158
    10 = 298
    11 = 355
    12 = 323
    13 = 721
    14 = 165
    15 = 370
    16 = 1320
    17 = 180
    18 = 679
    19 = 613
    20 = 255
    21 = 539
    22 = 801
    23 = 515
}

foreach ($row in $pocetValues.Keys) {
    $wsPocet.Range("AG$row").Value = $pocetValues[$row]
}

$wsPocet.Range("A24").Value = "Život během pandemie, Duševní zdraví, velikost dotázaného souboru celkově a ve skupinách, aktualizace 6. 10. 2021"
